# Apply the BOM workbook update:
# - Update "Design Last Modified" date
# - Update CONN-H16 unit price
# - Update two line-item totals in the Bill of Materials table
# - Update sheet view scroll position (topLeftCell) and window size

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Update the "Design Last Modified" date string (cell D8)
$ws.Range("D8").Value = "domingo, 13 de febrero de 2022"

# Update the unit price for CONN-H16 (cell F26)
$ws.Range("F26").Value = "€4,50"

# Update the recalculated totals for rows 18 and 19 (column J)
$ws.Range("J18").Value = 9.9499999999999993
$ws.Range("J19").Value = 22.2

# Scroll the sheet view so that row 16 is the top-left visible cell
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Application.ActiveWindow.ScrollColumn = 1

# Adjust the workbook window size recorded in the saved view
$excel.ActiveWindow.Width = 15345
$excel.ActiveWindow.Height = 6735
